$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update masthead text (volume number and week-covering date range) ---
$ws.Range("A8").Value = "Volume 29   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/24/2022  Through  10/30/2022"

# --- Plain numeric value updates ---
$ws.Range("M14").Value = -50
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -25
$ws.Range("N15").Value = -49.253731343283
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 33
$ws.Range("G16").Value = 34
$ws.Range("H16").Value = -2.941176470588
$ws.Range("I16").Value = 359
$ws.Range("J16").Value = 278
$ws.Range("K16").Value = 29.136690647482
$ws.Range("L16").Value = 37.022900763358
$ws.Range("M16").Value = 39.147286821705
$ws.Range("N16").Value = -66.100094428706
$ws.Range("C17").Value = 16
$ws.Range("D17").Value = 15
$ws.Range("E17").Value = 6.666666666666
$ws.Range("F17").Value = 47
$ws.Range("G17").Value = 61
$ws.Range("H17").Value = -22.950819672131
$ws.Range("I17").Value = 554
$ws.Range("J17").Value = 540
$ws.Range("K17").Value = 2.592592592592
$ws.Range("L17").Value = 22.566371681415
$ws.Range("M17").Value = 107.49063670412
$ws.Range("N17").Value = -32.848484848484
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -27.272727272727
$ws.Range("I18").Value = 277
$ws.Range("J18").Value = 169
$ws.Range("K18").Value = 63.905325443787
$ws.Range("L18").Value = 42.783505154639
$ws.Range("M18").Value = 118.110236220472
$ws.Range("N18").Value = -69.693654266958
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 46
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = -14.814814814814
$ws.Range("I19").Value = 434
$ws.Range("J19").Value = 393
$ws.Range("K19").Value = 10.432569974554
$ws.Range("L19").Value = 30.330330330330
$ws.Range("M19").Value = 114.851485148515
$ws.Range("N19").Value = 41.830065359477
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("G20").Value = 26
$ws.Range("H20").Value = -3.846153846153
$ws.Range("I20").Value = 251
$ws.Range("J20").Value = 182
$ws.Range("K20").Value = 37.912087912087
$ws.Range("L20").Value = 118.260869565217
$ws.Range("M20").Value = 164.210526315789
$ws.Range("N20").Value = -36.455696202531
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = 11.111111111111
$ws.Range("F21").Value = 170
$ws.Range("H21").Value = -15.422885572139
$ws.Range("I21").Value = 1914
$ws.Range("J21").Value = 1610
$ws.Range("K21").Value = 18.881987577639
$ws.Range("L21").Value = 37.697841726618
$ws.Range("M21").Value = 96.913580246913
$ws.Range("N21").Value = -46.729752296131
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 11
$ws.Range("K22").Value = 0
$ws.Range("C23").Value = 9
$ws.Range("D23").Value = 7
$ws.Range("E23").Value = 28.571428571428
$ws.Range("F23").Value = 33
$ws.Range("H23").Value = 26.923076923076
$ws.Range("I23").Value = 302
$ws.Range("J23").Value = 184
$ws.Range("K23").Value = 64.130434782608
$ws.Range("L23").Value = 105.442176870748
$ws.Range("M23").Value = 91.139240506329
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = 47.058823529411
$ws.Range("F24").Value = 84
$ws.Range("G24").Value = 71
$ws.Range("H24").Value = 18.309859154929
$ws.Range("I24").Value = 1056
$ws.Range("J24").Value = 791
$ws.Range("K24").Value = 33.501896333754
$ws.Range("L24").Value = 6.774519716885
$ws.Range("M24").Value = 63.975155279503
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 19
$ws.Range("E25").Value = -5.263157894736
$ws.Range("F25").Value = 53
$ws.Range("G25").Value = 85
$ws.Range("H25").Value = -37.647058823529
$ws.Range("I25").Value = 820
$ws.Range("J25").Value = 753
$ws.Range("K25").Value = 8.897742363877
$ws.Range("L25").Value = 15.492957746478
$ws.Range("M25").Value = 12.482853223594
$ws.Range("F26").Value = 6
$ws.Range("G26").Value = 6
$ws.Range("I26").Value = 53
$ws.Range("K26").Value = 8.163265306122
$ws.Range("L26").Value = 47.222222222222
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 9
$ws.Range("H27").Value = 28.571428571428
$ws.Range("I27").Value = 69
$ws.Range("J27").Value = 56
$ws.Range("K27").Value = 23.214285714285
$ws.Range("L27").Value = 72.5
$ws.Range("D28").Value = 1
$ws.Range("J28").Value = 66
$ws.Range("K28").Value = -51.515151515151
$ws.Range("M28").Value = -23.809523809523
$ws.Range("N28").Value = -65.957446808510
$ws.Range("D29").Value = 1
$ws.Range("J29").Value = 56
$ws.Range("K29").Value = -50
$ws.Range("M29").Value = -22.222222222222
$ws.Range("N29").Value = -69.892473118279

# --- Cells switching from numeric to text dash/placeholder style ---
# (value must be set with a leading apostrophe to force text type,
#  then formats are copied from a same-styled source cell so the
#  resulting style index matches a genuine "text" cell style.)
$ws.Range("D15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E15").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D26").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E26").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Cells switching from text dash/placeholder back to numeric ---
$ws.Range("D22").Value = 1
$ws.Range("G22").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E22").Value = -100
$ws.Range("K15").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
